$p = $ppt.ActivePresentation

# 1. Append a new "Title Only" slide (layout 11 = ppLayoutTitleOnly) at the
#    end of the deck (position 13) with the title "Other activities".
$newSlide = $p.Slides.Add($p.Slides.Count + 1, 11)
$titleRange = $newSlide.Shapes.Item(1).TextFrame.TextRange
$titleRange.Paragraphs(1, 1).Runs(1, 1).Text = "Other activities"
$titleRange.IndentLevel = 0
$titleRange.ParagraphFormat.Bullet.Visible = $false

# 2. Rename the (former) last slide's title from "External activities" to
#    "Other publications" (that slide keeps its existing two-column body).
$slide12 = $p.Slides.Item(12)
$slide12.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1, 1).Runs(1, 1).Text = "Other publications"

# 3. Reword the intro sentence on the "Data Platform" slide (slide 8) to
#    the past-tense / completed phrasing, leaving the remaining bullet
#    paragraphs untouched.
$slide8 = $p.Slides.Item(8)
$introRange = $slide8.Shapes.Item(2).TextFrame.TextRange
$introRange.Paragraphs(1, 1).Runs(1, 1).Text = "Within PNRR Agritech – Spoke 3, a Data Platform fostering collaboration and integration across research projects has been implemented."
